$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 4122.760395925205
$ws.Range("C3").Value = 4122.760395925205
$ws.Range("C4").Value = 3860.215426692186
$ws.Range("C5").Value = 3860.215426692186
$ws.Range("C6").Value = 3860.215426692186
$ws.Range("C7").Value = 3860.215426692186
$ws.Range("C8").Value = 3860.215426692186
$ws.Range("C9").Value = 3860.215426692186
$ws.Range("C10").Value = 3860.215426692186
$ws.Range("C11").Value = 3860.215426692186
$ws.Range("C12").Value = 3860.215426692186
